$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "RG567345646546"
$ws.Range("E3").Value = "MC345346574"
$ws.Range("E4").Value = "ZZZ344356546435"
